$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 708; existing rows 708:733 shift down to 709:734.
$ws.Rows.Item(708).Insert()

# Populate the newly inserted row 708 with the new record.
$ws.Cells.Item(708, 1).Value = 4
$ws.Cells.Item(708, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(708, 3).Value = 'Los Lagos'
$ws.Cells.Item(708, 4).Value = 45075
$ws.Cells.Item(708, 5).Value = 10
$ws.Cells.Item(708, 6).Value = 100112006
$ws.Cells.Item(708, 7).Value = 'Repollo'
$ws.Cells.Item(708, 8).Value = 'Crespo record'
$ws.Cells.Item(708, 9).Value = 'Segunda'
$ws.Cells.Item(708, 10).Value = 250
$ws.Cells.Item(708, 11).Value = 1500
$ws.Cells.Item(708, 12).Value = 1500
$ws.Cells.Item(708, 13).Value = 1500
$ws.Cells.Item(708, 14).Value = '$/unidad'
$ws.Cells.Item(708, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(708, 16).Value = 1500
$ws.Cells.Item(708, 17).Value = 1
$ws.Cells.Item(708, 18).Value = 'Hortaliza'
